# gui_single_search.py uses NAV excel (pretend DB) to get RFID from
# Dock#, DockLine#, PackageLine# inputs.
#
# This inserts a new "Unit ID" column (C) into the Dock/Line/RFID lookup
# sheet (shifting "Pkg Line No" -> D and "RFIDNumber" -> E) and refreshes
# the data rows to the latest NAV export order/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("Unit ID"). This shifts the existing
# "Pkg Line No" column from C->D and "RFIDNumber" from D->E.
$ws.Columns.Item(3).Insert()

# Header row - copy the bold/bordered header formatting onto the new
# header cell, then set its text.
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)
$ws.Cells.Item(1, 3).Value = "Unit ID"

# Final data set (Doc No, Doc Line No, Pkg Line No, RFIDNumber) in the
# new row order; the Unit ID column (C) is left blank for every row.
$data = @(
    @("S253441", 1, 3, "E2004702ED6060268CB9010D"),
    @("S243415", 2, 3, "E2004704D9906026AB7C010D"),
    @("S342451", 5, 8, "E2004704D9B06026AB7E0109"),
    @("S452341", 2, 1, "E20047053EB06026B1CE010A"),
    @("S412354", 4, 7, "E20047053E906026B1CC0105"),
    @("S234145", 1, 8, "E20047053EC06026B1CF0108"),
    @("S534241", 5, 5, "E20047053EA06026B1CD010A"),
    @("S253441", 1, 2, "E2004703EC9060269CAC0110"),
    @("S454132", 7, 3, "E2004704D9C06026AB7F0114"),
    @("S542314", 3, 1, "E2004704D9A06026AB7D010E")
)

$destRow = 2
foreach ($row in $data) {
    $ws.Cells.Item($destRow, 1).Value = $row[0]
    $ws.Cells.Item($destRow, 2).Value = $row[1]
    $ws.Cells.Item($destRow, 4).Value = $row[2]
    $ws.Cells.Item($destRow, 5).Value = $row[3]
    $destRow++
}

# Touch column C on each data row so an (empty) cell is materialized
# there, matching the source sheet's new blank "Unit ID" column.
$ws.Cells.Item(2, 1).Copy()
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
}
